# Updated cryptos list on Fri Jun 28 08:18:19 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'61.299.52"
$ws.Range("E2").Value = "  +0.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.427.14"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'574.33"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'145.39"
$ws.Range("E6").Value = "  +6.95%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "'3.426.53"
$ws.Range("E7").Value = "  +1.56%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  +2.03%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'7.65"
$ws.Range("E10").Value = "  +1.07%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +3.17%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "  +1.90%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'4.012.00"
$ws.Range("E13").Value = "  +1.52%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'28.04"
$ws.Range("E14").Value = "  +8.00%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.77%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.73%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'3.427.43"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'61.383.38"
$ws.Range("E18").Value = "  +0.84%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.27"
$ws.Range("E19").Value = "  +7.87%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'14.17"
$ws.Range("E20").Value = "  +3.14%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'9.39"
$ws.Range("E21").Value = "  +1.74%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'396.30"
$ws.Range("E22").Value = "  +6.55%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.565"
$ws.Range("E23").Value = "  +3.28%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'72.90"
$ws.Range("E24").Value = "  +3.07%  "

# Row 25 - Dai
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.46%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  -0.30%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -0.68%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "'3.576.56"
$ws.Range("E28").Value = "  +1.92%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  +1.14%  "

# Row 30 - RenderToken
$ws.Range("D30").Value = "'7.58"
$ws.Range("E30").Value = "  +3.39%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32 / Row 33 - Fetch.AI and InternetComputer(DFINITY) swap places
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.16"
$ws.Range("E32").Value = "  +2.08%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.47"
$ws.Range("E33").Value = "  -7.39%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +2.09%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  -0.08%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "'23.98"
$ws.Range("E36").Value = "  +3.14%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  +3.79%  "

# Row 38 - RenzoRestakedETH
$ws.Range("D38").Value = "'3.455.85"
$ws.Range("E38").Value = "  +1.81%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +1.82%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value = "  +0.41%  "

# Row 41 - Monero
$ws.Range("D41").Value = "'167.34"
$ws.Range("E41").Value = "  +1.55%  "

# Row 42 - Hedera
$ws.Range("D42").Value = "'0.0785"
$ws.Range("E42").Value = "  +3.67%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "'26.83"
$ws.Range("E43").Value = "  +8.12%  "

# Row 44 - Mantle
$ws.Range("D44").Value = "'0.798"
$ws.Range("E44").Value = "  +3.45%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  +0.74%  "

# Row 47 - Filecoin
$ws.Range("E47").Value = "  +3.92%  "

# Row 48 - OKB
$ws.Range("E48").Value = "  +0.50%  "

# Row 49 / Row 50 - Maker and ONDO swap places
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.16"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'2.577.53"
$ws.Range("E50").Value = "  +2.09%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  +2.56%  "
